{"js": "// The underlying commit only re-serialized word/document.xml: every hunk in\n// the diff reorders XML attributes within existing elements (namespace\n// declarations on <w:document>, the w:color on the \"self\" field run, the\n// sectPr page size/margins, the rPrDefault fonts/lang, the latentStyles\n// block, and the style/tblPr definitions). No text, formatting, structure,\n// or property VALUE changes anything - every attribute=value pair present\n// before is still present after, just written in a different order.\n//\n// Attribute emission order is a low-level OOXML-serializer concern; it is\n// not part of the Word document object model exposed by Office.js (there is\n// no API to reorder XML attributes), so there is no content-level edit to\n// perform here. We simply touch the body (a harmless, read-only load) so\n// the script still runs a real Office.js operation against the document\n// without altering any visible text, run formatting, or section/style\n// properties.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying commit only re-serialized word/document.xml: every hunk in\n# the diff reorders XML attributes within existing elements (namespace\n# declarations on <w:document>, the w:color on the \"self\" field run, the\n# sectPr page size/margins, the rPrDefault fonts/lang, the latentStyles\n# block, and the style/tblPr definitions). No text, formatting, structure,\n# or property VALUE changes anything - every attribute=value pair present\n# before is still present after, just written in a different order.\n#\n# Attribute emission order is a low-level OOXML-serializer concern; it is\n# not part of the Word COM object model (there is no property that controls\n# XML attribute order), so there is no content-level edit to perform here.\n# We simply touch the document (a harmless, read-only reference) so the\n# script still runs a real COM operation against the document without\n# altering any visible text, run formatting, or section/style properties.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
